# Update "江西-漫展信息.xlsx" to the state output at commit 456a3b4.
#
# Summary of changes:
#  - Sheet "展览"   (index 1): 16 "want-to-go count" / price refreshes, plus a
#    new con ("南昌·花绒万兽首届兽聚") inserted as row 34, pushing the two
#    trailing rows down by one (35, 36). Dimension A1:I35 -> A1:I36.
#  - Sheet "演出"   (index 2): 1 "want-to-go count" refresh.
#  - Sheet "本地生活" (index 3): untouched.
#  - Sheet "全部类型" (index 4): same refreshes as sheet 1 + the new Kpop row
#    that already lived only in this combined sheet, offset by one row versus
#    sheet 1 because of that extra "演出" row. Dimension A1:I36 -> A1:I37.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F3").Value  = 1849
$ws1.Range("F5").Value  = 5
$ws1.Range("F6").Value  = 814
$ws1.Range("G6").Value  = 36.6
$ws1.Range("F16").Value = 4349
$ws1.Range("F19").Value = 476
$ws1.Range("F20").Value = 416
$ws1.Range("F21").Value = 4
$ws1.Range("F22").Value = 989
$ws1.Range("F23").Value = 1740
$ws1.Range("F24").Value = 366
$ws1.Range("F26").Value = 14
$ws1.Range("F28").Value = 2044
$ws1.Range("F29").Value = 69
$ws1.Range("F31").Value = 5
$ws1.Range("F32").Value = 142
$ws1.Range("F33").Value = 92

# Insert a brand-new row at 34: old row 34 (代号鸢盛花行only) becomes row 35,
# old row 35 (ETI动漫节) becomes row 36.
$ws1.Rows.Item(34).Insert()

# Pick up the bordered/bold numbering style used by column A from the row
# above so the new index cell matches its neighbours.
$ws1.Range("A33").Copy($ws1.Range("A34"))

$ws1.Range("A34").Value = 33
$ws1.Range("B34").Value = "'2024-05-18"
$ws1.Range("C34").Value = "南昌·花绒万兽首届兽聚"
$ws1.Range("D34").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$ws1.Range("E34").Value = "2024.05.18 09:30-05.19 16:30"
$ws1.Range("F34").Value = 21
$ws1.Range("G34").Value = 50
$ws1.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=83689"
$ws1.Range("I34").Value = "//i2.hdslb.com/bfs/openplatform/202403/h4iL6IvI1711790121140.jpeg"

# Row that used to be 34 is now 35: refresh its sequential index + its one
# changed field (want-to-go count 210 -> 211); everything else shifted as-is.
$ws1.Range("A35").Value = 34
$ws1.Range("F35").Value = 211

# Row that used to be 35 is now 36: only the sequential index needs restating
# (content is identical to before the shift).
$ws1.Range("A36").Value = 35

# ---------------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 6

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 -- no changes
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (same data as sheet 1, offset by the extra 演出 row that
# lives at row 16 here)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F3").Value  = 1849
$ws4.Range("F5").Value  = 5
$ws4.Range("F6").Value  = 814
$ws4.Range("F16").Value = 6
$ws4.Range("F17").Value = 4349
$ws4.Range("F20").Value = 476
$ws4.Range("F21").Value = 416
$ws4.Range("F22").Value = 4
$ws4.Range("F23").Value = 989
$ws4.Range("F24").Value = 1740
$ws4.Range("F25").Value = 366
$ws4.Range("F27").Value = 14
$ws4.Range("F29").Value = 2044
$ws4.Range("F30").Value = 69
$ws4.Range("F32").Value = 5
$ws4.Range("F33").Value = 143
$ws4.Range("F34").Value = 92

# Insert the matching new row at 35 (one below where it landed in sheet 1).
$ws4.Rows.Item(35).Insert()

$ws4.Range("A34").Copy($ws4.Range("A35"))

$ws4.Range("A35").Value = 34
$ws4.Range("B35").Value = "'2024-05-18"
$ws4.Range("C35").Value = "南昌·花绒万兽首届兽聚"
$ws4.Range("D35").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$ws4.Range("E35").Value = "2024.05.18 09:30-05.19 16:30"
$ws4.Range("F35").Value = 21
$ws4.Range("G35").Value = 50
$ws4.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=83689"
$ws4.Range("I35").Value = "//i2.hdslb.com/bfs/openplatform/202403/h4iL6IvI1711790121140.jpeg"

$ws4.Range("A36").Value = 35
$ws4.Range("F36").Value = 211

$ws4.Range("A37").Value = 36

Write-Host "Applied 江西-漫展信息 refresh (456a3b4)."
